$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: year 2021 data
$ws.Range("K4").Value = 2021
$ws.Range("K5").Value = 7999.5
$ws.Range("K6").Value = $null
$ws.Range("K7").Formula = "=K5-K8"
$ws.Range("K8").Value = 252.9
$ws.Range("K9").Value = $null
$ws.Range("K10").Value = 690.4
$ws.Range("K11").Value = 968.2
$ws.Range("K12").Value = 655
$ws.Range("K13").Value = 691.2
$ws.Range("K14").Value = 1248.5
$ws.Range("K15").Value = 959.1
$ws.Range("K16").Value = 2596.6
$ws.Range("K17").Value = 133.6
$ws.Range("K18").Value = 57

# Copy styles from column J to column K so formatting matches
$ws.Range("J4:J18").Copy()
$ws.Range("K4:K18").PasteSpecial(-4122)  # xlPasteFormats

# Update selection to match the target view state
$ws.Range("N20").Select()
